$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156 (shifts existing rows 156-254 down to 157-255,
# and the used range grows to A1:R255).
$ws.Rows(156).Insert()

# Populate the newly inserted row 156 with the new weekly data point.
$ws.Cells.Item(156, 1).Value = 8
$ws.Cells.Item(156, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 44879
$ws.Cells.Item(156, 5).Value = 4
$ws.Cells.Item(156, 6).Value = 100112037
$ws.Cells.Item(156, 7).Value = "Cebollín"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 1000
$ws.Cells.Item(156, 11).Value = 1200
$ws.Cells.Item(156, 12).Value = 1400
$ws.Cells.Item(156, 13).Value = 1300
$ws.Cells.Item(156, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(156, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(156, 16).Value = 217
$ws.Cells.Item(156, 17).Value = 6
$ws.Cells.Item(156, 18).Value = "Hortaliza"
